# Update automatico via Actualizar 05-22-2020 07-29-23
#
# Adds the new daily record (2020-05-21 -> serial 43972) as row 70 of the
# "Condicion_Pacientes" table, extending the table/autofilter range from
# A1:F69 to A1:F70, then reflects the cell that was last being edited in
# Excel's UI (C69) as the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row - this is what Excel does natively when a new
# record is appended right below an Excel Table; it keeps the table/
# autoFilter ref, the sheet dimension and the styling consistent.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $lo.ListRows.Add()

# Copy the formatting (date format on column A, centered numbers on B:F)
# from the previous last row so the new row's cells reuse the existing
# styles instead of minting new ones.
$ws.Range("A69:F69").Copy()
$ws.Range("A70:F70").PasteSpecial(-4122)

# New day's figures.
$ws.Range("A70").Value = 43972
$ws.Range("B70").Value = 505
$ws.Range("C70").Value = 104
$ws.Range("D70").Value = 235
$ws.Range("E70").Value = 14
$ws.Range("F70").Value = 17

# Match the workbook's final on-screen selection.
$ws.Range("C69").Select()
